$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 2.5
$ws.Range("H4").Value = 2.75
$ws.Range("I4").Value = 3.3
$ws.Range("J4").Value = 3.5
$ws.Range("K4").Value = 1.83
$ws.Range("L4").Value = 4
$ws.Range("M4").Value = 1.17
$ws.Range("N4").Value = 5
$ws.Range("O4").Value = 1.62
$ws.Range("P4").Value = 2.2
$ws.Range("Q4").Value = 3.1
$ws.Range("R4").Value = 1.36
$ws.Range("X4").Value = 10
$ws.Range("Z4").Value = 26
$ws.Range("AA4").Value = 29
$ws.Range("AF4").Value = 81
$ws.Range("AI4").Value = 13
$ws.Range("AJ4").Value = 34
$ws.Range("AK4").Value = 34
$ws.Range("AN4").Value = 4.33
$ws.Range("AO4").Value = 17
$ws.Range("AX4").Value = 21
$ws.Range("BA4").Value = 126
$ws.Range("BB4").Value = 500

# Row 7
$ws.Range("G7").Value = 2.15
$ws.Range("I7").Value = 4.1
$ws.Range("J7").Value = 3
$ws.Range("L7").Value = 5
$ws.Range("M7").Value = 1.17
$ws.Range("N7").Value = 5
$ws.Range("X7").Value = 8
$ws.Range("Z7").Value = 19
$ws.Range("AA7").Value = 23
$ws.Range("AD7").Value = 6
$ws.Range("AG7").Value = 7.5
$ws.Range("AH7").Value = 19
$ws.Range("AI7").Value = 17
$ws.Range("AN7").Value = 3.75
$ws.Range("AO7").Value = 13

# Row 11
$ws.Range("I11").Value = 1.78
$ws.Range("J11").Value = 4.9
$ws.Range("L11").Value = 2.37
$ws.Range("W11").Value = 10
$ws.Range("Y11").Value = 15.5
$ws.Range("AA11").Value = 55
$ws.Range("AB11").Value = 65
$ws.Range("AC11").Value = 7.4
$ws.Range("AG11").Value = 5.6
$ws.Range("AL11").Value = 35
$ws.Range("AN11").Value = 6.1
$ws.Range("AO11").Value = 27
$ws.Range("AQ11").Value = 175
$ws.Range("AR11").Value = 250
$ws.Range("AU11").Value = 7.8
$ws.Range("AW11").Value = 3.45
$ws.Range("AX11").Value = 8.75
$ws.Range("BA11").Value = 70

# Row 12
$ws.Range("G12").Value = 1.72
$ws.Range("H12").Value = 3.8
$ws.Range("I12").Value = 4.1
$ws.Range("J12").Value = 2.25
$ws.Range("K12").Value = 2.27
$ws.Range("L12").Value = 4.3
$ws.Range("U12").Value = 1.62
$ws.Range("X12").Value = 9
$ws.Range("Z12").Value = 14
$ws.Range("AA12").Value = 13
$ws.Range("AD12").Value = 7.5
$ws.Range("AE12").Value = 14
$ws.Range("AF12").Value = 55
$ws.Range("AH12").Value = 26
$ws.Range("AN12").Value = 3.7
$ws.Range("AO12").Value = 8.25
$ws.Range("AP12").Value = 15.5
$ws.Range("AQ12").Value = 26
$ws.Range("AT12").Value = 3.05
$ws.Range("AX12").Value = 22
$ws.Range("BB12").Value = 300

# Row 33
$ws.Range("G33").Value = 2.85
$ws.Range("H33").Value = 3.5
$ws.Range("I33").Value = 2.18
$ws.Range("J33").Value = 3.4
$ws.Range("L33").Value = 2.75
$ws.Range("N33").Value = 7.8
$ws.Range("O33").Value = 1.27
$ws.Range("P33").Value = 3.45
$ws.Range("Q33").Value = 1.8
$ws.Range("R33").Value = 1.93
$ws.Range("T33").Value = 2.92
$ws.Range("U33").Value = 1.7
$ws.Range("V33").Value = 2.05
$ws.Range("W33").Value = 10
$ws.Range("X33").Value = 15
$ws.Range("Y33").Value = 10.5
$ws.Range("AA33").Value = 23
$ws.Range("AC33").Value = 7.8
$ws.Range("AD33").Value = 6.9
$ws.Range("AF33").Value = 60
$ws.Range("AJ33").Value = 21
$ws.Range("AK33").Value = 17
$ws.Range("AL33").Value = 26
$ws.Range("AN33").Value = 4.9
$ws.Range("AO33").Value = 15
$ws.Range("AQ33").Value = 65
$ws.Range("AT33").Value = 2.92
$ws.Range("AW33").Value = 4.2
$ws.Range("AX33").Value = 11

# Row 34
$ws.Range("H34").Value = 3.2
$ws.Range("I34").Value = 2.77
$ws.Range("J34").Value = 2.92
$ws.Range("K34").Value = 2.12
$ws.Range("N34").Value = 7
$ws.Range("Q34").Value = 1.98
$ws.Range("R34").Value = 1.78
$ws.Range("T34").Value = 2.75
$ws.Range("W34").Value = 8
$ws.Range("X34").Value = 11.75
$ws.Range("AC34").Value = 7
$ws.Range("AD34").Value = 6.3
$ws.Range("AG34").Value = 8.5
$ws.Range("AH34").Value = 14
$ws.Range("AJ34").Value = 32
$ws.Range("AP34").Value = 19.5
$ws.Range("AT34").Value = 2.75
$ws.Range("AU34").Value = 7
$ws.Range("BA34").Value = 110

# Row 37
$ws.Range("G37").Value = 2.27
$ws.Range("H37").Value = 3.1
$ws.Range("I37").Value = 3.05
$ws.Range("J37").Value = 2.85
$ws.Range("K37").Value = 2.05
$ws.Range("L37").Value = 3.6
$ws.Range("M37").Value = 1.06
$ws.Range("N37").Value = 8.77
$ws.Range("P37").Value = 2.62
$ws.Range("U37").Value = 1.83
$ws.Range("V37").Value = 1.78
$ws.Range("W37").Value = 6.9
$ws.Range("X37").Value = 10.5
$ws.Range("Y37").Value = 9.25
$ws.Range("Z37").Value = 23
$ws.Range("AB37").Value = 32
$ws.Range("AC37").Value = 7.9
$ws.Range("AD37").Value = 6
$ws.Range("AE37").Value = 15.5
$ws.Range("AF37").Value = 80
$ws.Range("AH37").Value = 15
$ws.Range("AI37").Value = 11.25
$ws.Range("AJ37").Value = 40
$ws.Range("AM37").Value = 800
$ws.Range("AN37").Value = 4.05
$ws.Range("AO37").Value = 11.75
$ws.Range("AP37").Value = 20
$ws.Range("AU37").Value = 7.1
$ws.Range("AV37").Value = 65
$ws.Range("AW37").Value = 4.85
$ws.Range("AX37").Value = 17
$ws.Range("AY37").Value = 25
$ws.Range("AZ37").Value = 80
$ws.Range("BB37").Value = 350

